$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new price entry as row 13.
$ws.Range("A13").Value = "HD SSD Kingston SA400S37 480GB"
$ws.Range("B13").Value = "259,98"
$ws.Range("C13").Value = "08/08/2025 - 22:48:13"
$ws.Range("D13").Value = "https://www.amazon.com.br/HD-SSD-KINGSTON-SA400S37-480GB/dp/B075BKXSCQ/ref=sr_1_3?crid=3ABFVIMZS2O0T&dib=eyJ2IjoiMSJ9.xtLqN2YY2lmra89PEqy7G0Y84YkphyyUWV5twVivkNse1ODjkVlE7xfYt0-FIma8U9i0nv0Se2nMSPd-hyWG2Teo__6mPZ9JRy8ISaS7yBPhkUYWYTAJopOQ4hRMxZ7dzQoHYb3lI3LHfw_YRtgQrv4Fwxhs4tWHKz4EMr7VTkiSDRSWgvk3N6BZT1FtUAJMQc6JscsuzwuHCoTXPfmZm36OgM4cVt-aZFb9XNsALKZBTBNEIumTM7NwH3bVuyT-z9NAVFNpSltZOaTQnUckcRMTscf3YKI2hCTLpSqWQig.WagU7fhpIkWdmNgtfShtJccPMhTWlh9y46wTOG-nbZY&dib_tag=se&keywords=ssd+500gb&qid=1732729870&sprefix=ssd%2Caps%2C160&sr=8-3&ufe=app_do%3Aamzn1.fos.6a09f7ec-d911-4889-ad70-de8dd83c8a74"
